$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Update AI.mean (Q) and AI.sd (R) columns per revised calculations
$ws.Range("Q3").Value = 0.42
$ws.Range("R3").Value = 0.86
$ws.Range("Q4").Value = 0.76
$ws.Range("R4").Value = 0.97
$ws.Range("Q5").Value = 0.95
$ws.Range("R5").Value = 1.41
$ws.Range("Q7").Value = 0.89
$ws.Range("R7").Value = 1.22
$ws.Range("Q8").Value = 0.47
$ws.Range("R8").Value = 0.74
$ws.Range("R9").Value = 1.1
$ws.Range("Q10").Value = 0.78
$ws.Range("R10").Value = 1
$ws.Range("R11").Value = 0.75
$ws.Range("Q12").Value = 0.68
$ws.Range("R12").Value = 0.97
$ws.Range("Q13").Value = 0.38
$ws.Range("R13").Value = 0.76
$ws.Range("Q14").Value = 0.97
$ws.Range("R14").Value = 1.19
$ws.Range("Q15").Value = 0.38
$ws.Range("R15").Value = 0.65
$ws.Range("Q17").Value = 0.68
$ws.Range("R17").Value = 0.97
$ws.Range("Q25").Value = 0.39
$ws.Range("R25").Value = 0.69
$ws.Range("Q43").Value = 0.37
$ws.Range("R43").Value = 0.66
$ws.Range("Q44").Value = 0.47
$ws.Range("R44").Value = 0.71
$ws.Range("Q46").Value = 0.77
$ws.Range("R46").Value = 1.11
$ws.Range("Q47").Value = 0.67
$ws.Range("R47").Value = 0.87
$ws.Range("Q48").Value = 1.06
$ws.Range("R48").Value = 1.13
$ws.Range("Q49").Value = 0.83
$ws.Range("R49").Value = 0.96
$ws.Range("Q50").Value = 0.51
$ws.Range("Q51").Value = 1.05
$ws.Range("R51").Value = 1.34
$ws.Range("Q52").Value = 0.82
$ws.Range("R52").Value = 0.93
$ws.Range("Q53").Value = 0.8
$ws.Range("R53").Value = 1.08
$ws.Range("R54").Value = 1.01
